# Disability Labour Force Participation widget.
#
# The "Data" sheet had an extra spacer column (D) removed, shifting the
# jurisdiction columns (NSW..Aust) one slot to the left, and the workbook's
# active tab moved from "Description" back to "Data".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Remove the blank spacer column D - this shifts E:O left to D:N and
# automatically updates the sheet dimension and the A1:O1 merged header
# down to A1:N1.
$ws.Columns.Item(4).Delete() | Out-Null

# The wrapped header / year rows reflow slightly now that the merged title
# cell is a column narrower.
$ws.Rows.Item(1).RowHeight = 29.85
$ws.Rows.Item(4).RowHeight = 15.65
$ws.Rows.Item(9).RowHeight = 15.65

# Make "Data" the active sheet/tab again (it had been "Description"),
# with C6 selected.
$ws.Activate() | Out-Null
$ws.Range("C6").Select() | Out-Null
